$win = $excel.ActiveWindow
$win.Width = 20490
$win.Height = 7755
Write-Host "done"
